$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells (I1 = "I0", J1 = "IF"), matching the style of the
# existing header cells (e.g. H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in data rows 2-28: column I is always 1, column J equals column H
for ($r = 2; $r -le 28; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 8).Value()
}
